# TDS-Data.xlsx edit script
# Commit: Data as on 9-1-21-Morning

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-empty calibration data cells (rows 35-100, 111-122) ---
$ws.Cells.Item(35, 3).Value = 0.078974000000000003
$ws.Cells.Item(35, 6).Value = 22
$ws.Cells.Item(35, 7).Value = 56
$ws.Cells.Item(35, 8).Value = 46

$ws.Cells.Item(36, 3).Value = 0.074138999999999997
$ws.Cells.Item(36, 6).Value = 21
$ws.Cells.Item(36, 7).Value = 49
$ws.Cells.Item(36, 8).Value = 45

$ws.Cells.Item(37, 3).Value = 0.064468999999999999
$ws.Cells.Item(37, 6).Value = 19
$ws.Cells.Item(37, 7).Value = 46
$ws.Cells.Item(37, 8).Value = 45

$ws.Cells.Item(38, 3).Value = 0.056410000000000002
$ws.Cells.Item(38, 6).Value = 18
$ws.Cells.Item(38, 7).Value = 45
$ws.Cells.Item(38, 8).Value = 47

$ws.Cells.Item(39, 3).Value = 0.051575000000000003
$ws.Cells.Item(39, 6).Value = 17
$ws.Cells.Item(39, 7).Value = 45
$ws.Cells.Item(39, 8).Value = 45

$ws.Cells.Item(40, 3).Value = 0.047545999999999998
$ws.Cells.Item(40, 6).Value = 16
$ws.Cells.Item(40, 7).Value = 47
$ws.Cells.Item(40, 8).Value = 45

$ws.Cells.Item(41, 3).Value = 0.039487000000000001
$ws.Cells.Item(41, 6).Value = 14
$ws.Cells.Item(41, 7).Value = 45
$ws.Cells.Item(41, 8).Value = 46

$ws.Cells.Item(42, 3).Value = 0.03304
$ws.Cells.Item(42, 6).Value = 12
$ws.Cells.Item(42, 7).Value = 45
$ws.Cells.Item(42, 8).Value = 45

$ws.Cells.Item(43, 3).Value = 0.026592999999999999
$ws.Cells.Item(43, 6).Value = 10
$ws.Cells.Item(43, 7).Value = 44
$ws.Cells.Item(43, 8).Value = 39

$ws.Cells.Item(44, 3).Value = 0.022564000000000001
$ws.Cells.Item(44, 6).Value = 9
$ws.Cells.Item(44, 7).Value = 45
$ws.Cells.Item(44, 8).Value = 42

$ws.Cells.Item(45, 3).Value = 0.015311
$ws.Cells.Item(45, 6).Value = 7
$ws.Cells.Item(45, 7).Value = 40
$ws.Cells.Item(45, 8).Value = 44

$ws.Cells.Item(46, 3).Value = 0.058021999999999997
$ws.Cells.Item(46, 6).Value = 16
$ws.Cells.Item(46, 7).Value = 40
$ws.Cells.Item(46, 8).Value = 46

$ws.Cells.Item(47, 3).Value = 0.053186999999999998
$ws.Cells.Item(47, 6).Value = 15
$ws.Cells.Item(47, 7).Value = 43
$ws.Cells.Item(47, 8).Value = 39

$ws.Cells.Item(48, 3).Value = 0.047545999999999998
$ws.Cells.Item(48, 6).Value = 14
$ws.Cells.Item(48, 7).Value = 45
$ws.Cells.Item(48, 8).Value = 40

$ws.Cells.Item(49, 3).Value = 0.038681
$ws.Cells.Item(49, 6).Value = 12
$ws.Cells.Item(49, 7).Value = 44
$ws.Cells.Item(49, 8).Value = 45

$ws.Cells.Item(50, 3).Value = 0.037874999999999999
$ws.Cells.Item(50, 6).Value = 12
$ws.Cells.Item(50, 7).Value = 40
$ws.Cells.Item(50, 8).Value = 44

$ws.Cells.Item(51, 3).Value = 0.033846000000000001
$ws.Cells.Item(51, 6).Value = 11
$ws.Cells.Item(51, 7).Value = 43
$ws.Cells.Item(51, 8).Value = 40

$ws.Cells.Item(52, 3).Value = 0.024982000000000001
$ws.Cells.Item(52, 6).Value = 9
$ws.Cells.Item(52, 7).Value = 40
$ws.Cells.Item(52, 8).Value = 41

$ws.Cells.Item(53, 3).Value = 0.018534999999999999
$ws.Cells.Item(53, 6).Value = 7
$ws.Cells.Item(53, 7).Value = 41
$ws.Cells.Item(53, 8).Value = 41

$ws.Cells.Item(54, 3).Value = 0.0137
$ws.Cells.Item(54, 6).Value = 5
$ws.Cells.Item(54, 7).Value = 38
$ws.Cells.Item(54, 8).Value = 40

$ws.Cells.Item(55, 3).Value = 0.011282
$ws.Cells.Item(55, 6).Value = 5
$ws.Cells.Item(55, 7).Value = 38
$ws.Cells.Item(55, 8).Value = 40

$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 35
$ws.Cells.Item(56, 8).Value = 39

$ws.Cells.Item(57, 3).Value = 0.051575000000000003
$ws.Cells.Item(57, 6).Value = 15
$ws.Cells.Item(57, 7).Value = 39
$ws.Cells.Item(57, 8).Value = 46

$ws.Cells.Item(58, 3).Value = 0.045934000000000003
$ws.Cells.Item(58, 6).Value = 13
$ws.Cells.Item(58, 7).Value = 43
$ws.Cells.Item(58, 8).Value = 39

$ws.Cells.Item(59, 3).Value = 0.038681
$ws.Cells.Item(59, 6).Value = 12
$ws.Cells.Item(59, 7).Value = 39
$ws.Cells.Item(59, 8).Value = 43

$ws.Cells.Item(60, 3).Value = 0.03304
$ws.Cells.Item(60, 6).Value = 10
$ws.Cells.Item(60, 7).Value = 37
$ws.Cells.Item(60, 8).Value = 43

$ws.Cells.Item(61, 3).Value = 0.028205000000000001
$ws.Cells.Item(61, 6).Value = 9
$ws.Cells.Item(61, 7).Value = 38
$ws.Cells.Item(61, 8).Value = 41

$ws.Cells.Item(62, 3).Value = 0.024982000000000001
$ws.Cells.Item(62, 6).Value = 8
$ws.Cells.Item(62, 7).Value = 38
$ws.Cells.Item(62, 8).Value = 39

$ws.Cells.Item(63, 3).Value = 0.016923000000000001
$ws.Cells.Item(63, 6).Value = 6
$ws.Cells.Item(63, 7).Value = 38
$ws.Cells.Item(63, 8).Value = 39

$ws.Cells.Item(64, 3).Value = 0.012893999999999999
$ws.Cells.Item(64, 6).Value = 5
$ws.Cells.Item(64, 7).Value = 38
$ws.Cells.Item(64, 8).Value = 38

$ws.Cells.Item(65, 3).Value = 0.0056410000000000002
$ws.Cells.Item(65, 6).Value = 2
$ws.Cells.Item(65, 7).Value = 37
$ws.Cells.Item(65, 8).Value = 38

$ws.Cells.Item(66, 3).Value = 0.0016119999999999999
$ws.Cells.Item(66, 6).Value = 1
$ws.Cells.Item(66, 7).Value = 35
$ws.Cells.Item(66, 8).Value = 38

$ws.Cells.Item(67, 3).Value = 0
$ws.Cells.Item(67, 4).Value = 0
$ws.Cells.Item(67, 5).Value = 0
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 31
$ws.Cells.Item(67, 8).Value = 36

$ws.Cells.Item(68, 3).Value = 0.029817
$ws.Cells.Item(68, 6).Value = 8
$ws.Cells.Item(68, 7).Value = 42
$ws.Cells.Item(68, 8).Value = 33

$ws.Cells.Item(69, 3).Value = 0.024176
$ws.Cells.Item(69, 6).Value = 7
$ws.Cells.Item(69, 7).Value = 33
$ws.Cells.Item(69, 8).Value = 39

$ws.Cells.Item(70, 3).Value = 0.017729000000000002
$ws.Cells.Item(70, 6).Value = 5
$ws.Cells.Item(70, 7).Value = 31
$ws.Cells.Item(70, 8).Value = 35

$ws.Cells.Item(71, 3).Value = 0.011282
$ws.Cells.Item(71, 6).Value = 4
$ws.Cells.Item(71, 7).Value = 31
$ws.Cells.Item(71, 8).Value = 34

$ws.Cells.Item(72, 3).Value = 0.0080590000000000002
$ws.Cells.Item(72, 6).Value = 3
$ws.Cells.Item(72, 7).Value = 31
$ws.Cells.Item(72, 8).Value = 34

$ws.Cells.Item(73, 3).Value = 0.0040289999999999996
$ws.Cells.Item(73, 6).Value = 1
$ws.Cells.Item(73, 7).Value = 31
$ws.Cells.Item(73, 8).Value = 33

$ws.Cells.Item(74, 3).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 31
$ws.Cells.Item(74, 8).Value = 33

$ws.Cells.Item(75, 3).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 30
$ws.Cells.Item(75, 8).Value = 31

$ws.Cells.Item(76, 3).Value = 0
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 29
$ws.Cells.Item(76, 8).Value = 30

$ws.Cells.Item(77, 3).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 7).Value = 30
$ws.Cells.Item(77, 8).Value = 27

$ws.Cells.Item(78, 3).Value = 0
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 25
$ws.Cells.Item(78, 8).Value = 28

$ws.Cells.Item(79, 3).Value = 0.014505000000000001
$ws.Cells.Item(79, 6).Value = 4
$ws.Cells.Item(79, 7).Value = 37
$ws.Cells.Item(79, 8).Value = 30

$ws.Cells.Item(80, 3).Value = 0.012893999999999999
$ws.Cells.Item(80, 6).Value = 4
$ws.Cells.Item(80, 7).Value = 34
$ws.Cells.Item(80, 8).Value = 29

$ws.Cells.Item(81, 3).Value = 0.0088640000000000004
$ws.Cells.Item(81, 6).Value = 3
$ws.Cells.Item(81, 7).Value = 29
$ws.Cells.Item(81, 8).Value = 32

$ws.Cells.Item(82, 3).Value = 0.040259000000000003
$ws.Cells.Item(82, 6).Value = 4
$ws.Cells.Item(82, 7).Value = 28
$ws.Cells.Item(82, 8).Value = 30

$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 27
$ws.Cells.Item(83, 8).Value = 30

$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 28
$ws.Cells.Item(84, 8).Value = 29

$ws.Cells.Item(85, 3).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 27
$ws.Cells.Item(85, 8).Value = 28

$ws.Cells.Item(86, 3).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 27
$ws.Cells.Item(86, 8).Value = 28

$ws.Cells.Item(87, 3).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 21
$ws.Cells.Item(87, 8).Value = 27

$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 24
$ws.Cells.Item(88, 8).Value = 25

$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(89, 4).Value = 0
$ws.Cells.Item(89, 5).Value = 0
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 24
$ws.Cells.Item(89, 8).Value = 25

$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(90, 4).Value = 0
$ws.Cells.Item(90, 5).Value = 0
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 31
$ws.Cells.Item(90, 8).Value = 29

$ws.Cells.Item(91, 3).Value = 0
$ws.Cells.Item(91, 4).Value = 0
$ws.Cells.Item(91, 5).Value = 0
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 29
$ws.Cells.Item(91, 8).Value = 29

$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 28
$ws.Cells.Item(92, 8).Value = 25

$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 28
$ws.Cells.Item(93, 8).Value = 26

$ws.Cells.Item(94, 3).Value = 0
$ws.Cells.Item(94, 4).Value = 0
$ws.Cells.Item(94, 5).Value = 0
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 7).Value = 26
$ws.Cells.Item(94, 8).Value = 26

$ws.Cells.Item(95, 3).Value = 0
$ws.Cells.Item(95, 4).Value = 0
$ws.Cells.Item(95, 5).Value = 0
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 25
$ws.Cells.Item(95, 8).Value = 26

$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 24
$ws.Cells.Item(96, 8).Value = 25

$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 23
$ws.Cells.Item(97, 8).Value = 24

$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 23
$ws.Cells.Item(98, 8).Value = 20

$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 19
$ws.Cells.Item(99, 8).Value = 23

$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 4).Value = 0
$ws.Cells.Item(100, 5).Value = 0
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 21
$ws.Cells.Item(100, 8).Value = 21

$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 0
$ws.Cells.Item(111, 5).Value = 0
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 0

$ws.Cells.Item(112, 3).Value = 1.452161
$ws.Cells.Item(112, 6).Value = 355
$ws.Cells.Item(112, 7).Value = 439
$ws.Cells.Item(112, 8).Value = 408

$ws.Cells.Item(113, 3).Value = 1.410256
$ws.Cells.Item(113, 6).Value = 355
$ws.Cells.Item(113, 7).Value = 409
$ws.Cells.Item(113, 8).Value = 414

$ws.Cells.Item(114, 3).Value = 1.349817
$ws.Cells.Item(114, 6).Value = 354
$ws.Cells.Item(114, 7).Value = 393
$ws.Cells.Item(114, 8).Value = 386

$ws.Cells.Item(115, 3).Value = 1.293407
$ws.Cells.Item(115, 6).Value = 354
$ws.Cells.Item(115, 7).Value = 379
$ws.Cells.Item(115, 8).Value = 384

$ws.Cells.Item(116, 3).Value = 1.2369300000000001
$ws.Cells.Item(116, 6).Value = 350
$ws.Cells.Item(116, 7).Value = 387
$ws.Cells.Item(116, 8).Value = 363

$ws.Cells.Item(117, 3).Value = 1.2023440000000001
$ws.Cells.Item(117, 6).Value = 351
$ws.Cells.Item(117, 7).Value = 395
$ws.Cells.Item(117, 8).Value = 357

$ws.Cells.Item(118, 3).Value = 1.1394869999999999
$ws.Cells.Item(118, 6).Value = 349
$ws.Cells.Item(118, 7).Value = 385
$ws.Cells.Item(118, 8).Value = 337

$ws.Cells.Item(119, 3).Value = 1.070989
$ws.Cells.Item(119, 6).Value = 345
$ws.Cells.Item(119, 7).Value = 382
$ws.Cells.Item(119, 8).Value = 333

$ws.Cells.Item(120, 3).Value = 1.033919
$ws.Cells.Item(120, 6).Value = 345
$ws.Cells.Item(120, 7).Value = 324
$ws.Cells.Item(120, 8).Value = 376

$ws.Cells.Item(121, 3).Value = 0.99201499999999998
$ws.Cells.Item(121, 6).Value = 344
$ws.Cells.Item(121, 7).Value = 311
$ws.Cells.Item(121, 8).Value = 372

$ws.Cells.Item(122, 3).Value = 0.93640999999999996
$ws.Cells.Item(122, 6).Value = 344
$ws.Cells.Item(122, 7).Value = 362
$ws.Cells.Item(122, 8).Value = 302

# --- Row 364: new data point (H364), I364 now holds its own literal formula ---
$ws.Cells.Item(364, 8).Value = 12
$ws.Range("I364").Formula = "=(H364+G364)/2"

# --- N7, N9, N10: replaced with manually edited (non-shared) formulas ---
$ws.Range("N7").Formula = "=(126*100*J7+332*100*K7+70*L7*100)/(100*J7+100*K7+L7*100+M7*100)"
$ws.Range("N9").Formula = "=(126*100*J9+324*100*K9+70*L9*100)/(100*J9+100*K9+L9*100+M9*100)"
$ws.Range("N10").Formula = "=(126*100*J10+324*100*K10+70*L10*100)/(100*J10+100*K10+L10*100+M10*100)"

# --- Update the active cell selection to match the saved view state ---
[void]$ws.Range("N10").Select()
